$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Pdgfa"
$ws.Range("C2").Value = "Pdgfra"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 13.08425366666667
$ws.Range("H2").Value = 39.252761
$ws.Range("I2").Value = 0.4321946987699228
$ws.Range("J2").Value = 0.4321946987699228
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 1.155747666666667
$ws.Range("N2").Value = 3.467243
$ws.Range("O2").Value = 0.004246591903937912
$ws.Range("P2").Value = 0.004246591903937912
$ws.Range("Q2").Value = 15.12209564532478
$ws.Range("R2").Value = 136.098860807923
$ws.Range("S2").Value = 0.001835354508721239
$ws.Range("T2").Value = 0.001835354508721239

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Pdgfa"
$ws.Range("C3").Value = "Pdgfra"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 13.08425366666667
$ws.Range("H3").Value = 39.252761
$ws.Range("I3").Value = 0.4321946987699228
$ws.Range("J3").Value = 0.4321946987699228
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 270.7963256666667
$ws.Range("N3").Value = 812.3889770000001
$ws.Range("O3").Value = 0.9949935590256014
$ws.Range("P3").Value = 0.9949935590256014
$ws.Range("Q3").Value = 3543.167817023945
$ws.Range("R3").Value = 31888.5103532155
$ws.Range("S3").Value = 0.4300309415210832
$ws.Range("T3").Value = 0.4300309415210832

# Row 4
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Pdgfa"
$ws.Range("C4").Value = "Pdgfra"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 13.08425366666667
$ws.Range("H4").Value = 39.252761
$ws.Range("I4").Value = 0.4321946987699228
$ws.Range("J4").Value = 0.4321946987699228
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.2067996666666667
$ws.Range("N4").Value = 0.620399
$ws.Range("O4").Value = 0.0007598490704606447
$ws.Range("P4").Value = 0.0007598490704606446
$ws.Range("Q4").Value = 2.705819296848778
$ws.Range("R4").Value = 24.352373671639
$ws.Range("S4").Value = 0.0003284027401183442
$ws.Range("T4").Value = 0.0003284027401183442

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Pdgfa"
$ws.Range("C5").Value = "Pdgfra"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 0.7043496666666668
$ws.Range("H5").Value = 2.113049
$ws.Range("I5").Value = 0.0232658430330821
$ws.Range("J5").Value = 0.0232658430330821
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 1.155747666666667
$ws.Range("N5").Value = 3.467243
$ws.Range("O5").Value = 0.004246591903937912
$ws.Range("P5").Value = 0.004246591903937912
$ws.Range("Q5").Value = 0.8140504837674445
$ws.Range("R5").Value = 7.326454353907001
$ws.Range("S5").Value = 0.00009880054066257672
$ws.Range("T5").Value = 0.00009880054066257672

# Row 6
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Pdgfa"
$ws.Range("C6").Value = "Pdgfra"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 0.7043496666666668
$ws.Range("H6").Value = 2.113049
$ws.Range("I6").Value = 0.0232658430330821
$ws.Range("J6").Value = 0.0232658430330821
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 270.7963256666667
$ws.Range("N6").Value = 812.3889770000001
$ws.Range("O6").Value = 0.9949935590256014
$ws.Range("P6").Value = 0.9949935590256014
$ws.Range("Q6").Value = 190.7353017178748
$ws.Range("R6").Value = 1716.617715460873
$ws.Range("S6").Value = 0.02314936396321735
$ws.Range("T6").Value = 0.02314936396321735

# Row 7
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Pdgfa"
$ws.Range("C7").Value = "Pdgfra"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 0.7043496666666668
$ws.Range("H7").Value = 2.113049
$ws.Range("I7").Value = 0.0232658430330821
$ws.Range("J7").Value = 0.0232658430330821
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 0.2067996666666667
$ws.Range("N7").Value = 0.620399
$ws.Range("O7").Value = 0.0007598490704606447
$ws.Range("P7").Value = 0.0007598490704606446
$ws.Range("Q7").Value = 0.1456592762834445
$ws.Range("R7").Value = 1.310933486551
$ws.Range("S7").Value = 0.0000176785292021707
$ws.Range("T7").Value = 0.0000176785292021707

# Row 8
$ws.Range("A8").Value = "sCs"
$ws.Range("B8").Value = "Pdgfa"
$ws.Range("C8").Value = "Pdgfra"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 16.48537666666667
$ws.Range("H8").Value = 49.45613
$ws.Range("I8").Value = 0.5445394581969951
$ws.Range("J8").Value = 0.5445394581969951
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 0.6666666666666666
$ws.Range("M8").Value = 1.155747666666667
$ws.Range("N8").Value = 3.467243
$ws.Range("O8").Value = 0.004246591903937912
$ws.Range("P8").Value = 0.004246591903937912
$ws.Range("Q8").Value = 19.05293561662111
$ws.Range("R8").Value = 171.47642054959
$ws.Range("S8").Value = 0.002312436854554097
$ws.Range("T8").Value = 0.002312436854554097

# Row 9
$ws.Range("A9").Value = "sCs"
$ws.Range("B9").Value = "Pdgfa"
$ws.Range("C9").Value = "Pdgfra"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 16.48537666666667
$ws.Range("H9").Value = 49.45613
$ws.Range("I9").Value = 0.5445394581969951
$ws.Range("J9").Value = 0.5445394581969951
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 270.7963256666667
$ws.Range("N9").Value = 812.3889770000001
$ws.Range("O9").Value = 0.9949935590256014
$ws.Range("P9").Value = 0.9949935590256014
$ws.Range("Q9").Value = 4464.179428564335
$ws.Range("R9").Value = 40177.61485707902
$ws.Range("S9").Value = 0.5418132535413008
$ws.Range("T9").Value = 0.5418132535413008

# Row 10
$ws.Range("A10").Value = "sCs"
$ws.Range("B10").Value = "Pdgfa"
$ws.Range("C10").Value = "Pdgfra"
$ws.Range("D10").Value = "sCs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 16.48537666666667
$ws.Range("H10").Value = 49.45613
$ws.Range("I10").Value = 0.5445394581969951
$ws.Range("J10").Value = 0.5445394581969951
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 0.2067996666666667
$ws.Range("N10").Value = 0.620399
$ws.Range("O10").Value = 0.0007598490704606447
$ws.Range("P10").Value = 0.0007598490704606446
$ws.Range("Q10").Value = 3.409170399541112
$ws.Range("R10").Value = 30.68253359587
$ws.Range("S10").Value = 0.0004137678011401299
$ws.Range("T10").Value = 0.0004137678011401298
